# Update wcl_weights (column Q) from 6000 to 2000 for the rows where it
# currently holds the placeholder value of 6000 between rows 97 and 171.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 97..171) {
    $cell = $ws.Cells.Item($r, 17)  # Column Q = 17
    if ($cell.Value2 -eq 6000) {
        $cell.Value2 = 2000
    }
}
